# Insert a new row at row 174 that duplicates the "Z / YARD BLOCK" row
# (currently row 2: Section=Z, Block Number=151, Block Length=50,
# Speed Limit=100, Infrastructure=YARD BLOCK) just before the existing
# "SWITCH TO YARD" row, pushing every row from 174 downward down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 174 - this shifts rows 174..223 down to 175..224
# and copies formatting from the row above, same as Excel's own
# Insert Row behaviour.
$ws.Rows("174:174").Insert()

# Populate the new row 174 with the same values/styles as row 2 (the
# template "YARD BLOCK" entry) by copy/pasting that row onto it.
$ws.Range("A2:E2").Copy($ws.Range("A174:E174"))

# Reflect the resulting selection, matching what Excel shows right after
# performing this insert-row-and-fill operation.
$ws.Range("A174:E174").Select()
